$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 28.479168
$ws.Cells.Item(2, 8).Value = 85.437504
$ws.Cells.Item(2, 9).Value = 0.4446244458164738
$ws.Cells.Item(2, 10).Value = 0.4446244458164738
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 127.6999736666667
$ws.Cells.Item(2, 14).Value = 383.099921
$ws.Cells.Item(2, 15).Value = 0.9554352891750322
$ws.Cells.Item(2, 16).Value = 0.9554352891750322
$ws.Cells.Item(2, 17).Value = 3636.789003648576
$ws.Cells.Item(2, 18).Value = 32731.10103283718
$ws.Cells.Item(2, 19).Value = 0.4248098859629511
$ws.Cells.Item(2, 20).Value = 0.4248098859629511

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 28.479168
$ws.Cells.Item(3, 8).Value = 85.437504
$ws.Cells.Item(3, 9).Value = 0.4446244458164738
$ws.Cells.Item(3, 10).Value = 0.4446244458164738
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.4321196666666667
$ws.Cells.Item(3, 14).Value = 1.296359
$ws.Cells.Item(3, 15).Value = 0.003233065495828321
$ws.Cells.Item(3, 16).Value = 0.003233065495828321
$ws.Cells.Item(3, 17).Value = 12.306408583104
$ws.Cells.Item(3, 18).Value = 110.757677247936
$ws.Cells.Item(3, 19).Value = 0.00143749995437103
$ws.Cells.Item(3, 20).Value = 0.00143749995437103

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 28.479168
$ws.Cells.Item(4, 8).Value = 85.437504
$ws.Cells.Item(4, 9).Value = 0.4446244458164738
$ws.Cells.Item(4, 10).Value = 0.4446244458164738
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.77305
$ws.Cells.Item(4, 14).Value = 14.31915
$ws.Cells.Item(4, 15).Value = 0.03571136528892854
$ws.Cells.Item(4, 16).Value = 0.03571136528892854
$ws.Cells.Item(4, 17).Value = 135.9324928224
$ws.Cells.Item(4, 18).Value = 1223.3924354016
$ws.Cells.Item(4, 19).Value = 0.01587814600093951
$ws.Cells.Item(4, 20).Value = 0.01587814600093951

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 28.479168
$ws.Cells.Item(5, 8).Value = 85.437504
$ws.Cells.Item(5, 9).Value = 0.4446244458164738
$ws.Cells.Item(5, 10).Value = 0.4446244458164738
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.751186
$ws.Cells.Item(5, 14).Value = 2.253558
$ws.Cells.Item(5, 15).Value = 0.00562028004021099
$ws.Cells.Item(5, 16).Value = 0.00562028004021099
$ws.Cells.Item(5, 17).Value = 21.393152293248
$ws.Cells.Item(5, 18).Value = 192.538370639232
$ws.Cells.Item(5, 19).Value = 0.002498913898212201
$ws.Cells.Item(5, 20).Value = 0.002498913898212201

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 18.12667766666667
$ws.Cells.Item(6, 8).Value = 54.380033
$ws.Cells.Item(6, 9).Value = 0.2829985767855128
$ws.Cells.Item(6, 10).Value = 0.2829985767855128
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 127.6999736666667
$ws.Cells.Item(6, 14).Value = 383.099921
$ws.Cells.Item(6, 15).Value = 0.9554352891750322
$ws.Cells.Item(6, 16).Value = 0.9554352891750322
$ws.Cells.Item(6, 17).Value = 2314.776260697488
$ws.Cells.Item(6, 18).Value = 20832.98634627739
$ws.Cells.Item(6, 19).Value = 0.270386827047189
$ws.Cells.Item(6, 20).Value = 0.270386827047189

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 18.12667766666667
$ws.Cells.Item(7, 8).Value = 54.380033
$ws.Cells.Item(7, 9).Value = 0.2829985767855128
$ws.Cells.Item(7, 10).Value = 0.2829985767855128
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.4321196666666667
$ws.Cells.Item(7, 14).Value = 1.296359
$ws.Cells.Item(7, 15).Value = 0.003233065495828321
$ws.Cells.Item(7, 16).Value = 0.003233065495828321
$ws.Cells.Item(7, 17).Value = 7.832893911094112
$ws.Cells.Item(7, 18).Value = 70.496045199847
$ws.Cells.Item(7, 19).Value = 0.0009149529339737631
$ws.Cells.Item(7, 20).Value = 0.0009149529339737631

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 18.12667766666667
$ws.Cells.Item(8, 8).Value = 54.380033
$ws.Cells.Item(8, 9).Value = 0.2829985767855128
$ws.Cells.Item(8, 10).Value = 0.2829985767855128
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 4.77305
$ws.Cells.Item(8, 14).Value = 14.31915
$ws.Cells.Item(8, 15).Value = 0.03571136528892854
$ws.Cells.Item(8, 16).Value = 0.03571136528892854
$ws.Cells.Item(8, 17).Value = 86.51953883688336
$ws.Cells.Item(8, 18).Value = 778.6758495319502
$ws.Cells.Item(8, 19).Value = 0.01010626555183434
$ws.Cells.Item(8, 20).Value = 0.01010626555183434

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 18.12667766666667
$ws.Cells.Item(9, 8).Value = 54.380033
$ws.Cells.Item(9, 9).Value = 0.2829985767855128
$ws.Cells.Item(9, 10).Value = 0.2829985767855128
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.751186
$ws.Cells.Item(9, 14).Value = 2.253558
$ws.Cells.Item(9, 15).Value = 0.00562028004021099
$ws.Cells.Item(9, 16).Value = 0.00562028004021099
$ws.Cells.Item(9, 17).Value = 13.61650648971267
$ws.Cells.Item(9, 18).Value = 122.548558407414
$ws.Cells.Item(9, 19).Value = 0.001590531252515735
$ws.Cells.Item(9, 20).Value = 0.001590531252515735

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.513346
$ws.Cells.Item(10, 8).Value = 34.540038
$ws.Cells.Item(10, 9).Value = 0.179749460544048
$ws.Cells.Item(10, 10).Value = 0.179749460544048
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 127.6999736666667
$ws.Cells.Item(10, 14).Value = 383.099921
$ws.Cells.Item(10, 15).Value = 0.9554352891750322
$ws.Cells.Item(10, 16).Value = 0.9554352891750322
$ws.Cells.Item(10, 17).Value = 1470.253981015222
$ws.Cells.Item(10, 18).Value = 13232.285829137
$ws.Cells.Item(10, 19).Value = 0.1717389778139586
$ws.Cells.Item(10, 20).Value = 0.1717389778139586

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 11.513346
$ws.Cells.Item(11, 8).Value = 34.540038
$ws.Cells.Item(11, 9).Value = 0.179749460544048
$ws.Cells.Item(11, 10).Value = 0.179749460544048
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.4321196666666667
$ws.Cells.Item(11, 14).Value = 1.296359
$ws.Cells.Item(11, 15).Value = 0.003233065495828321
$ws.Cells.Item(11, 16).Value = 0.003233065495828321
$ws.Cells.Item(11, 17).Value = 4.975143235738
$ws.Cells.Item(11, 18).Value = 44.776289121642
$ws.Cells.Item(11, 19).Value = 0.0005811417787787159
$ws.Cells.Item(11, 20).Value = 0.0005811417787787159

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 11.513346
$ws.Cells.Item(12, 8).Value = 34.540038
$ws.Cells.Item(12, 9).Value = 0.179749460544048
$ws.Cells.Item(12, 10).Value = 0.179749460544048
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 4.77305
$ws.Cells.Item(12, 14).Value = 14.31915
$ws.Cells.Item(12, 15).Value = 0.03571136528892854
$ws.Cells.Item(12, 16).Value = 0.03571136528892854
$ws.Cells.Item(12, 17).Value = 54.95377612530001
$ws.Cells.Item(12, 18).Value = 494.5839851277001
$ws.Cells.Item(12, 19).Value = 0.006419098645976346
$ws.Cells.Item(12, 20).Value = 0.006419098645976346

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 11.513346
$ws.Cells.Item(13, 8).Value = 34.540038
$ws.Cells.Item(13, 9).Value = 0.179749460544048
$ws.Cells.Item(13, 10).Value = 0.179749460544048
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.751186
$ws.Cells.Item(13, 14).Value = 2.253558
$ws.Cells.Item(13, 15).Value = 0.00562028004021099
$ws.Cells.Item(13, 16).Value = 0.00562028004021099
$ws.Cells.Item(13, 17).Value = 8.648664328356
$ws.Cells.Item(13, 18).Value = 77.83797895520401
$ws.Cells.Item(13, 19).Value = 0.001010242305334406
$ws.Cells.Item(13, 20).Value = 0.001010242305334406

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 5.932994999999999
$ws.Cells.Item(14, 8).Value = 17.798985
$ws.Cells.Item(14, 9).Value = 0.09262751685396531
$ws.Cells.Item(14, 10).Value = 0.09262751685396531
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 127.6999736666667
$ws.Cells.Item(14, 14).Value = 383.099921
$ws.Cells.Item(14, 15).Value = 0.9554352891750322
$ws.Cells.Item(14, 16).Value = 0.9554352891750322
$ws.Cells.Item(14, 17).Value = 757.6433052644649
$ws.Cells.Item(14, 18).Value = 6818.789747380184
$ws.Cells.Item(14, 19).Value = 0.08849959835093352
$ws.Cells.Item(14, 20).Value = 0.08849959835093352

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 5.932994999999999
$ws.Cells.Item(15, 8).Value = 17.798985
$ws.Cells.Item(15, 9).Value = 0.09262751685396531
$ws.Cells.Item(15, 10).Value = 0.09262751685396531
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.4321196666666667
$ws.Cells.Item(15, 14).Value = 1.296359
$ws.Cells.Item(15, 15).Value = 0.003233065495828321
$ws.Cells.Item(15, 16).Value = 0.003233065495828321
$ws.Cells.Item(15, 17).Value = 2.563763821735
$ws.Cells.Item(15, 18).Value = 23.073874395615
$ws.Cells.Item(15, 19).Value = 0.0002994708287048115
$ws.Cells.Item(15, 20).Value = 0.0002994708287048115

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 5.932994999999999
$ws.Cells.Item(16, 8).Value = 17.798985
$ws.Cells.Item(16, 9).Value = 0.09262751685396531
$ws.Cells.Item(16, 10).Value = 0.09262751685396531
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 4.77305
$ws.Cells.Item(16, 14).Value = 14.31915
$ws.Cells.Item(16, 15).Value = 0.03571136528892854
$ws.Cells.Item(16, 16).Value = 0.03571136528892854
$ws.Cells.Item(16, 17).Value = 28.31848178475
$ws.Cells.Item(16, 18).Value = 254.86633606275
$ws.Cells.Item(16, 19).Value = 0.00330785509017834
$ws.Cells.Item(16, 20).Value = 0.00330785509017834

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 5.932994999999999
$ws.Cells.Item(17, 8).Value = 17.798985
$ws.Cells.Item(17, 9).Value = 0.09262751685396531
$ws.Cells.Item(17, 10).Value = 0.09262751685396531
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.751186
$ws.Cells.Item(17, 14).Value = 2.253558
$ws.Cells.Item(17, 15).Value = 0.00562028004021099
$ws.Cells.Item(17, 16).Value = 0.00562028004021099
$ws.Cells.Item(17, 17).Value = 4.456782782069999
$ws.Cells.Item(17, 18).Value = 40.11104503863
$ws.Cells.Item(17, 19).Value = 0.0005205925841486482
$ws.Cells.Item(17, 20).Value = 0.0005205925841486482
